# Extend the "SanPossidonio" style sheet with new daily rows 270-301
# (data updated through 28/06 included), replicating the existing
# formatting of column A (date style) used by the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 269
$firstNewSerial = 44344   # 2021-05-28, day after the last existing row (44343)
$lastNewRow = 301

# Copy the formatting of the last existing row's date cell (A269) so the
# new date cells (A270:A301) keep the same style (s="2": centered, bordered,
# bold, date/time number format).
$ws.Range("A$lastRow").Copy()
$ws.Range("A" + ($lastRow + 1) + ":A" + $lastNewRow).PasteSpecial(-4122)

for ($row = $lastRow + 1; $row -le $lastNewRow; $row++) {
    $serial = $firstNewSerial + ($row - ($lastRow + 1))
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}
